# Lesson 7 ImageSheet-Module7 update.
# Fills in / corrects the image-name list in rows 31-41 of Sheet1 and
# moves the active selection to C40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B31").Value = "n/a"
$ws.Range("B32").Value = "n/a"
$ws.Range("B33").Value = "n/a"
$ws.Range("B34").Value = "DLC.png"
$ws.Range("B35").Value = "n/a"
$ws.Range("B36").Value = "mapPlanning.jpg"

$ws.Range("B37").Value = "image33.jpeg"
$ws.Range("C37").Clear()
$ws.Range("D37").Clear()

$ws.Range("B38").Value = "binoculars.jpg"
$ws.Range("B39").Value = "SBLTER_Metadata.png"

$ws.Range("B40").Value = "SBLTER_attirbutes.png"
$ws.Range("C40").Clear()

$ws.Range("B41").Value = "SBLTER_metadataFields.png"

$ws.Range("C40").Select()
